$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.838.61"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "3.149.00"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'574.10"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "'149.01"
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.148.55"

$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = "  -3.66%  "

$ws.Range("D11").Value = "'6.12"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").Value = "'36.94"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").Value = "3.663.22"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "64.960.33"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "3.146.35"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "'7.08"
$ws.Range("E18").Value = "  -1.55%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'503.70"
$ws.Range("E20").Value = "  -1.43%  "

$ws.Range("D21").Value = "'14.76"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").Value = "'0.711"
$ws.Range("E22").Value = "  -2.89%  "

$ws.Range("D23").Value = "'15.10"
$ws.Range("E23").Value = "  -2.75%  "

$ws.Range("D24").Value = "'7.68"
$ws.Range("E24").Value = "  -2.65%  "

$ws.Range("D25").Value = "'83.81"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").Value = "'8.82"
$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("D28").Value = "'2.88"
$ws.Range("E28").Value = "  -1.85%  "

$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("D30").Value = "'2.80"
$ws.Range("E30").Value = "  +5.10%  "

$ws.Range("D31").Value = "'27.47"
$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").Value = "'6.17"
$ws.Range("E34").Value = "  +1.60%  "

$ws.Range("D35").Value = "'6.44"
$ws.Range("E35").Value = "  -3.22%  "

$ws.Range("D36").Value = "'54.63"
$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("D37").Value = "'0.0893"
$ws.Range("E37").Value = "  +4.11%  "

$ws.Range("D38").Value = "'473.65"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("E39").Value = "  -2.32%  "

$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  -2.45%  "

$ws.Range("D41").Value = "'8.62"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("D42").Value = "3.009.79"
$ws.Range("E42").Value = "  -3.79%  "

$ws.Range("E43").Value = "  -3.75%  "

$ws.Range("D44").Value = "'0.281"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").Value = "'2.41"
$ws.Range("E45").Value = "  -1.38%  "

$ws.Range("D46").Value = "'28.10"
$ws.Range("E46").Value = "  -4.37%  "

$ws.Range("D47").Value = "0.0₃0575"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("E49").Value = "  -2.32%  "

$ws.Range("E50").Value = "  -3.80%  "

$ws.Range("D51").Value = "'33.35"
$ws.Range("E51").Value = "  +5.64%  "
